$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 69; $r++) {
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 12).Value = $kVal
}

$ws.Range("C83").Value = 158
